# Regla falsa (false position) table update
# f(x) = x^2 - 4, root search between xi=1 and xs=3
# Extends the table from 7 to 10 iterations (rows 2-11) with updated values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1.75
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = -0.9375
$ws.Range("F2").Value = -3
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 1.000001

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1.94736842105263
$ws.Range("C3").Value = 1.75
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = -0.20775623268698
$ws.Range("F3").Value = -0.9375
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 0.101351351351351

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1.98936170212766
$ws.Range("C4").Value = 1.94736842105263
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = -0.0424400181077416
$ws.Range("F4").Value = -0.20775623268698
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 0.0211089220377145

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 1.99786780383795
$ws.Range("C5").Value = 1.98936170212766
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = -0.008524238387713899
$ws.Range("F5").Value = -0.0424400181077416
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 0.0042575898635301

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1.99957337883959
$ws.Range("C6").Value = 1.99786780383795
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = -0.0017063026360237
$ws.Range("F6").Value = -0.008524238387713899
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 0.0008529694482265

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 1.99991466848707
$ws.Range("C7").Value = 1.99957337883959
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = -0.0003413187702436
$ws.Range("F7").Value = -0.0017063026360237
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = 0.0001706521047421

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 1.99998293340615
$ws.Range("C8").Value = 1.99991466848707
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = -(6.826608412824339 * [Math]::Pow(10, -5))
$ws.Range("F8").Value = -0.0003413187702436
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 3.4132750804102 * [Math]::Pow(10, -5)

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 1.99999658666958
$ws.Range("C9").Value = 1.99998293340615
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = -(1.36533100323355 * [Math]::Pow(10, -5))
$ws.Range("F9").Value = -(6.826608412824339 * [Math]::Pow(10, -5))
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 6.82664336504653 * [Math]::Pow(10, -6)

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 1.99999931733345
$ws.Range("C10").Value = 1.99999658666958
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = -(2.73066573486247 * [Math]::Pow(10, -6))
$ws.Range("F10").Value = -(1.36533100323355 * [Math]::Pow(10, -5))
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 1.36533240132136 * [Math]::Pow(10, -6)

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 1.99999986346667
$ws.Range("C11").Value = 1.99999931733345
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = -(5.46133295920015 * [Math]::Pow(10, -7))
$ws.Range("F11").Value = -(2.73066573486247 * [Math]::Pow(10, -6))
$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 2.73066629408884 * [Math]::Pow(10, -7)

